$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.771.11'
$ws.Range('E2').Value = '  -6.59%  '
$ws.Range('D3').Value = '1.700.08'
$ws.Range('E3').Value = '  -5.68%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -5.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5087'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -14.54%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2608'
$ws.Range('D8').ClearFormats()
$ws.Range('E9').Value = '  -5.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06146'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -9.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07327'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('D12').Value = '1.684.84'
$ws.Range('E12').Value = '  -6.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.444'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.38%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.931.30'
$ws.Range('E14').Value = '  -5.61%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5747'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -8.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008197'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -10.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.63'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -12.71%  '
$ws.Range('D18').Value = '26.834.95'
$ws.Range('E18').Value = '  -6.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.032'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -7.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.75'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '185.68'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -11.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.247'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -8.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.47'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.678'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1150'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -9.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.32'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.325'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05657'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -9.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.332'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -6.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.479'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.449'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.669'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.009'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5921'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.636'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01600'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.39%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.070.67'
$ws.Range('E40').Value = '  -5.45%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.903'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -7.59%  '
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.003'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.87'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').Value = '1.857.79'
$ws.Range('E45').Value = '  -5.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.60'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.35%  '
$ws.Range('E47').Value = '  -4.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.012'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.142'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4332'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05210'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.30%  '
